# Auto-generated edits applying Rafflesia_Profits.xlsx diff
# Updates currentAveragePrice* / Leve Price* / Leve Profit* columns (H-N)
# for specific leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5: Met a Sticky End
$ws.Range("H5").Value = 86.8
$ws.Range("I5").Value = 86.8
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 86.8
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 28.2

# Row 43: Growing Is Knowing
$ws.Range("H43").Value = 2302
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2302
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 2302
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -2440

# Row 118: Crafty Concoctions
$ws.Range("H118").Value = 375
$ws.Range("I118").Value = 375
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 1125
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 532

# Row 127: Liquid Competence
$ws.Range("H127").Value = 375
$ws.Range("I127").Value = 375
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 1125
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 3835

# Row 131: Mindful Study
$ws.Range("H131").Value = 12500
$ws.Range("I131").Value = 5000
$ws.Range("J131").Value = 20000
$ws.Range("K131").Value = 15000
$ws.Range("L131").Value = 60000
$ws.Range("M131").Value = -9960

# Row 135: For Tired Minds
$ws.Range("H135").Value = 2000
$ws.Range("I135").Value = 1666.6666
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 14999.9994
$ws.Range("L135").Value = 27000
$ws.Range("M135").Value = -12464.9994

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 5000
$ws.Range("I137").Value = 5000
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 15000
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -12450

$ws = $wb.Worksheets.Item("ARM")
# Row 96: The Gauntlet Is Cast
$ws.Range("H96").Value = 34668.8
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 34668.8
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 34668.8
$ws.Range("N96").Value = -40160.8

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 387.75
$ws.Range("I86").Value = 387.75
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 387.75
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 735.25

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 387.75
$ws.Range("I89").Value = 387.75
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 1938.75
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 3677.25

# Row 92: Have Blade, Will Travel
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# Row 100: And My Axe
$ws.Range("H100").Value = 31698
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 31698
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 31698
$ws.Range("N100").Value = -33862

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 1995
$ws.Range("I105").Value = 1990
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1990
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -243
$ws.Range("N105").Value = -5494

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 7635.0835
$ws.Range("I134").Value = 4302.2
$ws.Range("J134").Value = 10015.714
$ws.Range("K134").Value = 12906.6
$ws.Range("L134").Value = 30047.142
$ws.Range("M134").Value = -10371.6

$ws = $wb.Worksheets.Item("CRP")
# Row 50: The Arsenal of Theocracy
$ws.Range("H50").Value = 19999.5
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 19999.5
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 19999.5
$ws.Range("N50").Value = -21249.5

# Row 51: Greenstone for Greenhorns
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 1999
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 1999
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 1999
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -2405

# Row 59: Bow Down to Magic
$ws.Range("H59").Value = 50000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 50000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 50000
$ws.Range("N59").Value = -52290

# Row 61: Incant Now, Think Later
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# Row 93: Reeling for Rods
$ws.Range("H93").Value = 26000
$ws.Range("I93").Value = 26000
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 26000
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -24128

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 3139
$ws.Range("I132").Value = 1699
$ws.Range("J132").Value = 3499
$ws.Range("K132").Value = 5097
$ws.Range("L132").Value = 10497
$ws.Range("M132").Value = -2567

# Row 136: Turali Quality
$ws.Range("H136").Value = 1999
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 1999
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 5997
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water
$ws.Range("H4").Value = 255.22223
$ws.Range("I4").Value = 255.22223
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 765.66669
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -653.66669

# Row 21: Shy Is the Oyster
$ws.Range("H21").Value = 78.40000000000001
$ws.Range("I21").Value = 66.333336
$ws.Range("J21").Value = 96.5
$ws.Range("K21").Value = 199.000008
$ws.Range("L21").Value = 289.5
$ws.Range("M21").Value = -26.00000800000001
$ws.Range("N21").Value = -635.5

# Row 38: Pretty as a Picture
$ws.Range("H38").Value = 1503.1666
$ws.Range("I38").Value = 1914.5
$ws.Range("J38").Value = 1297.5
$ws.Range("K38").Value = 5743.5
$ws.Range("L38").Value = 3892.5
$ws.Range("M38").Value = -5396.5
$ws.Range("N38").Value = -4586.5

# Row 60: Drinking to Your Health
$ws.Range("H60").Value = 360.57144
$ws.Range("I60").Value = 185
$ws.Range("J60").Value = 799.5
$ws.Range("K60").Value = 555
$ws.Range("L60").Value = 2398.5
$ws.Range("M60").Value = -304
$ws.Range("N60").Value = -2900.5

# Row 68: Such a Butter Face
$ws.Range("H68").Value = 778.875
$ws.Range("I68").Value = 701.5
$ws.Range("J68").Value = 856.25
$ws.Range("K68").Value = 2104.5
$ws.Range("L68").Value = 2568.75
$ws.Range("M68").Value = -1293.5
$ws.Range("N68").Value = -4190.75

# Row 71: No Margarine of Error (L)
$ws.Range("H71").Value = 778.875
$ws.Range("I71").Value = 701.5
$ws.Range("J71").Value = 856.25
$ws.Range("K71").Value = 6313.5
$ws.Range("L71").Value = 7706.25
$ws.Range("M71").Value = -2257.5
$ws.Range("N71").Value = -15818.25

# Row 97: The Frier Never Lies
$ws.Range("H97").Value = 421
$ws.Range("I97").Value = 92
$ws.Range("J97").Value = 750
$ws.Range("K97").Value = 276
$ws.Range("L97").Value = 2250
$ws.Range("M97").Value = 220

# Row 98: Sweet Kiss of Death
$ws.Range("H98").Value = 700
$ws.Range("I98").Value = 700
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2100
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -602
$ws.Range("N98").ClearContents()

# Row 107: Slippery Service
$ws.Range("H107").Value = 911.25
$ws.Range("I107").Value = 574.5
$ws.Range("J107").Value = 1248
$ws.Range("K107").Value = 1723.5
$ws.Range("L107").Value = 3744
$ws.Range("M107").Value = 196.5
$ws.Range("N107").Value = -7584

# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 1680.6
$ws.Range("I113").Value = 899
$ws.Range("J113").Value = 1876
$ws.Range("K113").Value = 2697
$ws.Range("L113").Value = 5628
$ws.Range("M113").Value = -527
$ws.Range("N113").Value = -9968

# Row 117: A Good Omen
$ws.Range("H117").Value = 5212.5
$ws.Range("I117").Value = 524.125
$ws.Range("J117").Value = 11463.667
$ws.Range("K117").Value = 1572.375
$ws.Range("L117").Value = 34391.001
$ws.Range("M117").Value = 1869.625
$ws.Range("N117").Value = -41275.001

# Row 121: A Cookie for Your Troubles
$ws.Range("H121").Value = 373.33334
$ws.Range("I121").Value = 62.5
$ws.Range("J121").Value = 995
$ws.Range("K121").Value = 187.5
$ws.Range("L121").Value = 2985
$ws.Range("M121").Value = 1122.5

# Row 122: Salt of the North
$ws.Range("H122").Value = 382.2
$ws.Range("I122").Value = 290.66666
$ws.Range("J122").Value = 519.5
$ws.Range("K122").Value = 2615.99994
$ws.Range("L122").Value = 4675.5
$ws.Range("M122").Value = -165.9999399999997
$ws.Range("N122").Value = -9575.5

# Row 123: Topping Up the Pot
$ws.Range("H123").Value = 4000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 4000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 12000
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -16900

# Row 129: Comfort Food
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("N129").ClearContents()

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 2222.3635
$ws.Range("I131").Value = 527.5
$ws.Range("J131").Value = 2599
$ws.Range("K131").Value = 1582.5
$ws.Range("L131").Value = 7797
$ws.Range("M131").Value = 3457.5
$ws.Range("N131").Value = -17877

# Row 132: More Mezcal
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 107: Whetstones for the Workers
$ws.Range("H107").Value = 1200
$ws.Range("I107").Value = 1200
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1200
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 720

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 3886.25
$ws.Range("I7").Value = 2966.818
$ws.Range("J7").Value = 14000
$ws.Range("K7").Value = 2966.818
$ws.Range("L7").Value = 14000
$ws.Range("M7").Value = -2854.818
$ws.Range("N7").Value = -14224

# Row 34: Breeches Served Cold
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 1610.75
$ws.Range("I93").Value = 1800.8572
$ws.Range("J93").Value = 280
$ws.Range("K93").Value = 1800.8572
$ws.Range("L93").Value = 280
$ws.Range("M93").Value = -552.8571999999999
$ws.Range("N93").Value = -2776

# Row 122: Hell on Leather
$ws.Range("H122").Value = 728
$ws.Range("I122").Value = 728
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2184
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 266

# Row 126: Battered Books
$ws.Range("H126").Value = 3886.25
$ws.Range("I126").Value = 2966.818
$ws.Range("J126").Value = 14000
$ws.Range("K126").Value = 8900.454000000002
$ws.Range("L126").Value = 42000
$ws.Range("M126").Value = -6430.454000000002
$ws.Range("N126").Value = -46940

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire
$ws.Range("H122").Value = 1216.5
$ws.Range("I122").Value = 796.1111
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 2388.3333
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = 61.66670000000022

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 2431.1875
$ws.Range("I136").Value = 2431.1875
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7293.5625
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4743.5625

